$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 7270
$ws.Range("I3").Value = 7491
$ws.Range("F4").Value = 1875
$ws.Range("I4").Value = 1720
$ws.Range("I5").Value = 708
$ws.Range("I6").Value = 8971
$ws.Range("F7").Value = 24064
$ws.Range("I7").Value = 26160

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I6").Value = 188
$ws.Range("I7").Value = 819
$ws.Range("I8").Value = 1544
$ws.Range("H11").Value = 345
$ws.Range("I14").Value = 142
$ws.Range("I15").Value = 299
$ws.Range("I18").Value = 207
$ws.Range("I19").Value = 734
$ws.Range("I23").Value = 252
$ws.Range("I29").Value = 1558
$ws.Range("I31").Value = 261
$ws.Range("I33").Value = 1147
$ws.Range("I36").Value = 360
$ws.Range("I37").Value = 806
$ws.Range("I42").Value = 1009
$ws.Range("I43").Value = 222
$ws.Range("I53").Value = 298
$ws.Range("I54").Value = 507
$ws.Range("I60").Value = 154
$ws.Range("F63").Value = 165
$ws.Range("H63").Value = 230
$ws.Range("I63").Value = 108
$ws.Range("I65").Value = 614
$ws.Range("I67").Value = 978
$ws.Range("I72").Value = 102
$ws.Range("I75").Value = 81
$ws.Range("I76").Value = 381
$ws.Range("I77").Value = 161
$ws.Range("I78").Value = 348
$ws.Range("I79").Value = 751
$ws.Range("I80").Value = 82
$ws.Range("I83").Value = 569
$ws.Range("I84").Value = 226
$ws.Range("I85").Value = 1161
$ws.Range("I88").Value = 244
$ws.Range("I89").Value = 308
$ws.Range("I95").Value = 404
$ws.Range("I96").Value = 307
$ws.Range("I98").Value = 189
$ws.Range("I99").Value = 448
$ws.Range("F101").Value = 24064
$ws.Range("I101").Value = 26160

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 326
$ws.Range("I3").Value = 441
$ws.Range("I6").Value = 304
$ws.Range("I7").Value = 1161

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("H4").Value = 26
$ws.Range("H7").Value = 345

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 457
$ws.Range("I3").Value = 448
$ws.Range("I7").Value = 1544

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I2").Value = 61
$ws.Range("I7").Value = 298

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I3").Value = 248
$ws.Range("I7").Value = 819

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("I3").Value = 71
$ws.Range("I7").Value = 308

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I6").Value = 126
$ws.Range("I7").Value = 307

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("I3").Value = 36
$ws.Range("I7").Value = 142

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 239
$ws.Range("I7").Value = 806

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I6").Value = 119
$ws.Range("I7").Value = 448

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I6").Value = 294
$ws.Range("I7").Value = 978

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I6").Value = 101
$ws.Range("I7").Value = 261

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("I2").Value = 83
$ws.Range("I6").Value = 55
$ws.Range("I7").Value = 226

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I6").Value = 190
$ws.Range("I7").Value = 614

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I2").Value = 190
$ws.Range("I7").Value = 569

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I3").Value = 137
$ws.Range("I5").Value = 22
$ws.Range("I7").Value = 404

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 268
$ws.Range("I3").Value = 416
$ws.Range("I6").Value = 366
$ws.Range("I7").Value = 1147

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I3").Value = 111
$ws.Range("I7").Value = 507

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 459
$ws.Range("I3").Value = 531
$ws.Range("I6").Value = 435
$ws.Range("I7").Value = 1558

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I3").Value = 213
$ws.Range("I5").Value = 23
$ws.Range("I6").Value = 234
$ws.Range("I7").Value = 734

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I4").Value = 40
$ws.Range("I6").Value = 180
$ws.Range("I7").Value = 381

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I3").Value = 47
$ws.Range("I6").Value = 55
$ws.Range("I7").Value = 188

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I3").Value = 280
$ws.Range("I6").Value = 428
$ws.Range("I7").Value = 1009

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I3").Value = 89
$ws.Range("I6").Value = 130
$ws.Range("I7").Value = 348

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I3").Value = 89
$ws.Range("I7").Value = 252

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 220
$ws.Range("I3").Value = 245
$ws.Range("I7").Value = 751

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("I6").Value = 99
$ws.Range("I7").Value = 207

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I3").Value = 118
$ws.Range("I6").Value = 113
$ws.Range("I7").Value = 360

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I2").Value = 90
$ws.Range("I7").Value = 299

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I3").Value = 17
$ws.Range("I7").Value = 189

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I3").Value = 79
$ws.Range("I7").Value = 244

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("I3").Value = 24
$ws.Range("I7").Value = 81

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("I2").Value = 55
$ws.Range("I7").Value = 154

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("I2").Value = 50
$ws.Range("I7").Value = 222

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("I6").Value = 51
$ws.Range("I7").Value = 102

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("I6").Value = 41
$ws.Range("I7").Value = 161

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("I6").Value = 47
$ws.Range("I7").Value = 82
